# Auto-generated script to apply scheduled market-data refresh to Sheets/Phantom_Profits.xlsx
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) per leve row,
# mirroring a scheduled runner re-pulling Universalis market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2612.4666
$ws.Range("J17").Value = 2735.5
$ws.Range("L17").Value = 8206.5
$ws.Range("N17").Value = -8542.5
$ws.Range("H19").Value = 1280.875
$ws.Range("I19").Value = 1274.5
$ws.Range("J19").Value = 1287.25
$ws.Range("K19").Value = 1274.5
$ws.Range("L19").Value = 1287.25
$ws.Range("M19").Value = -1099.5
$ws.Range("N19").Value = -1637.25
$ws.Range("H32").Value = 6247.5
$ws.Range("I32").Value = 3000.5
$ws.Range("K32").Value = 3000.5
$ws.Range("M32").Value = -2674.5
$ws.Range("H34").Value = 2862.5
$ws.Range("I34").Value = 2862.5
$ws.Range("K34").Value = 2862.5
$ws.Range("M34").Value = -2659.5
$ws.Range("H36").Value = 2862.5
$ws.Range("I36").Value = 2862.5
$ws.Range("K36").Value = 2862.5
$ws.Range("M36").Value = -2147.5
$ws.Range("H51").Value = 10849.8
$ws.Range("I51").Value = 11499.6
$ws.Range("K51").Value = 11499.6
$ws.Range("M51").Value = -11015.6
$ws.Range("H99").Value = 582
$ws.Range("I99").Value = 582
$ws.Range("K99").Value = 1746
$ws.Range("M99").Value = -248
$ws.Range("H103").Value = 529
$ws.Range("I103").Value = 293.5
$ws.Range("J103").Value = 1000
$ws.Range("K103").Value = 880.5
$ws.Range("L103").Value = 3000
$ws.Range("M103").Value = -294.5
$ws.Range("N103").Value = -4172
$ws.Range("H116").Value = 4340
$ws.Range("I116").Value = 2110.5
$ws.Range("K116").Value = 2110.5
$ws.Range("M116").Value = 1331.5
$ws.Range("H137").Value = 2860.1304
$ws.Range("J137").Value = 4727.1816
$ws.Range("L137").Value = 14181.5448
$ws.Range("N137").Value = -19281.5448
$ws.Range("H141").Value = 2573
$ws.Range("I141").Value = 2763.5386
$ws.Range("J141").Value = 96
$ws.Range("K141").Value = 8290.6158
$ws.Range("L141").Value = 288
$ws.Range("M141").Value = -3110.6158
$ws.Range("N141").Value = -10648

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2487.1516
$ws.Range("I32").Value = 1814.75
$ws.Range("K32").Value = 1814.75
$ws.Range("M32").Value = -1527.75
$ws.Range("H45").Value = 2374.7778
$ws.Range("I45").Value = 2374.7778
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2374.7778
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1997.7778
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 5170.2856
$ws.Range("J61").Value = 7140.8
$ws.Range("L61").Value = 7140.8
$ws.Range("N61").Value = -7564.8
$ws.Range("H97").Value = 1303.9445
$ws.Range("I97").Value = 792.41174
$ws.Range("K97").Value = 792.41174
$ws.Range("M97").Value = -296.41174
$ws.Range("H102").Value = 786.64703
$ws.Range("I102").Value = 786.64703
$ws.Range("K102").Value = 786.64703
$ws.Range("M102").Value = 835.35297
$ws.Range("H132").Value = 1555.8889
$ws.Range("I132").Value = 1555.8889
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4667.6667
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2137.6667
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 5170.2856
$ws.Range("J136").Value = 7140.8
$ws.Range("L136").Value = 21422.4
$ws.Range("N136").Value = -26522.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1108.8572
$ws.Range("I94").Value = 1219.625
$ws.Range("J94").Value = 961.1667
$ws.Range("K94").Value = 1219.625
$ws.Range("L94").Value = 961.1667
$ws.Range("M94").Value = -768.625
$ws.Range("N94").Value = -1863.1667
$ws.Range("H99").Value = 610.53845
$ws.Range("I99").Value = 575.0909
$ws.Range("K99").Value = 575.0909
$ws.Range("M99").Value = 922.9091
$ws.Range("H105").Value = 3892.25
$ws.Range("I105").Value = 3882.5454
$ws.Range("J105").Value = 3999
$ws.Range("K105").Value = 3882.5454
$ws.Range("L105").Value = 3999
$ws.Range("M105").Value = -2135.5454
$ws.Range("N105").Value = -7493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3077434.8
$ws.Range("I22").Value = 465.2
$ws.Range("K22").Value = 465.2
$ws.Range("M22").Value = -115.2
$ws.Range("H31").Value = 2839.8
$ws.Range("I31").Value = 2839.8
$ws.Range("K31").Value = 2839.8
$ws.Range("M31").Value = -2544.8
$ws.Range("H34").Value = 2839.8
$ws.Range("I34").Value = 2839.8
$ws.Range("K34").Value = 2839.8
$ws.Range("M34").Value = -2637.8
$ws.Range("H105").Value = 1075.4286
$ws.Range("I105").Value = 1054.6666
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 1054.6666
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = 692.3334
$ws.Range("N105").Value = -4694
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52080
$ws.Range("H132").Value = 1543.9474
$ws.Range("I132").Value = 1543.9474
$ws.Range("K132").Value = 4631.8422
$ws.Range("M132").Value = -2101.8422

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 125674
$ws.Range("I44").Value = 166724
$ws.Range("J44").Value = 2524
$ws.Range("K44").Value = 500172
$ws.Range("L44").Value = 7572
$ws.Range("M44").Value = -499774
$ws.Range("N44").Value = -8368
$ws.Range("H107").Value = 114.42857
$ws.Range("J107").Value = 116.6
$ws.Range("L107").Value = 349.8
$ws.Range("N107").Value = -4189.8
$ws.Range("H129").Value = 2908.3333
$ws.Range("J129").Value = 2277
$ws.Range("L129").Value = 6831
$ws.Range("N129").Value = -16831
$ws.Range("H131").Value = 2222.1667
$ws.Range("I131").Value = 2562.889
$ws.Range("J131").Value = 1200
$ws.Range("K131").Value = 7688.667
$ws.Range("L131").Value = 3600
$ws.Range("M131").Value = -2648.667
$ws.Range("N131").Value = -13680

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1980
$ws.Range("I132").Value = 1920
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 5760
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3230
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1748.75
$ws.Range("I46").Value = 1712.8572
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 1712.8572
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -1524.8572
$ws.Range("N46").Value = -2376
$ws.Range("H55").Value = 350.9
$ws.Range("I55").Value = 237.5
$ws.Range("J55").Value = 521
$ws.Range("K55").Value = 237.5
$ws.Range("L55").Value = 521
$ws.Range("M55").Value = -64.5
$ws.Range("N55").Value = -867
$ws.Range("H93").Value = 1216.2106
$ws.Range("I93").Value = 1357.2
$ws.Range("J93").Value = 687.5
$ws.Range("K93").Value = 1357.2
$ws.Range("L93").Value = 687.5
$ws.Range("M93").Value = -109.2
$ws.Range("N93").Value = -3183.5
$ws.Range("H100").Value = 1821.6
$ws.Range("I100").Value = 1752
$ws.Range("K100").Value = 1752
$ws.Range("M100").Value = -1211
$ws.Range("H132").Value = 2358.8
$ws.Range("I132").Value = 2321.4546
$ws.Range("K132").Value = 6964.3638
$ws.Range("M132").Value = -4434.3638
$ws.Range("H136").Value = 6572.467
$ws.Range("J136").Value = 10624.75
$ws.Range("L136").Value = 31874.25
$ws.Range("N136").Value = -36974.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1932.4375
$ws.Range("I100").Value = 652.1
$ws.Range("J100").Value = 4066.3333
$ws.Range("K100").Value = 1304.2
$ws.Range("L100").Value = 8132.6666
$ws.Range("M100").Value = -763.2
$ws.Range("N100").Value = -9214.6666
$ws.Range("H109").Value = 59722
$ws.Range("J109").Value = 59722
$ws.Range("L109").Value = 59722
$ws.Range("N109").Value = -62496
$ws.Range("H132").Value = 2517.2896
$ws.Range("I132").Value = 2268.2778
$ws.Range("K132").Value = 6804.8334
$ws.Range("M132").Value = -4274.8334
$ws.Range("H136").Value = 21537.375
$ws.Range("I136").Value = 22081.4
$ws.Range("J136").Value = 20630.666
$ws.Range("K136").Value = 66244.20000000001
$ws.Range("L136").Value = 61891.99800000001
$ws.Range("M136").Value = -63694.20000000001
$ws.Range("N136").Value = -66991.99800000001
